$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 681 (shifts existing rows 681-781 down to 682-782)
$ws.Rows(681).Insert()

# Populate the newly inserted row 681 with the new weekly record
$ws.Range("A681").Value = 6
$ws.Range("B681").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C681").Value = "Metropolitana"
$ws.Range("D681").Value = 45131
$ws.Range("E681").Value = 13
$ws.Range("F681").Value = 100112044
$ws.Range("G681").Value = "Perejil"
$ws.Range("H681").Value = "Sin especificar"
$ws.Range("I681").Value = "Primera"
$ws.Range("J681").Value = 260
$ws.Range("K681").Value = 14000
$ws.Range("L681").Value = 15000
$ws.Range("M681").Value = 14423
$ws.Range("N681").Value = '$/docena de atados'
$ws.Range("O681").Value = "Región Metropolitana"
$ws.Range("P681").Value = 4808
$ws.Range("Q681").Value = 3
$ws.Range("R681").Value = "Hortaliza"
